# Generate Report for Handback
#
# Updates the localization-status report to reflect that the de-de handback
# has now completed in sync with en-US:
#   - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#     (Overview zh-cn/de-de columns, and the Status column on each language
#     sheet)
#   - The "Latest Target File" / "Latest Handback File" / "Latest Handback
#     DateTime" columns on the zh-cn and de-de sheets are now populated with
#     the generated handback info, with the target-file cell turned into a
#     hyperlink back to the source markdown file (mirroring column A).
#   - A few status/report columns are widened so the new text fits.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdFile   = "30d6b8d0-65a7-4bca-8fda-8549ca162525.md"
$mdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/911d1209e298439f1c81bce9a831435666804619/e2e/$mdFile"

$zhXlf    = "30d6b8d0-65a7-4bca-8fda-8549ca162525.922bb0427acb37c44ea70b4f17270d17034ef070.zh-cn.xlf"
$deXlf    = "30d6b8d0-65a7-4bca-8fda-8549ca162525.922bb0427acb37c44ea70b4f17270d17034ef070.de-de.xlf"

$zhHandbackDate = "2016-08-24 11:05:47"
$deHandbackDate = "2016-08-24 11:05:54"

# --- Overview sheet: widen the status columns and refresh their text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

foreach ($cellRef in @("E2", "F2", "E3", "F3")) {
    if ($overview.Range($cellRef).Value() -eq $oldStatus) {
        $overview.Range($cellRef).Value = $newStatus
    }
}

# --- Per-language detail sheets ---
foreach ($langInfo in @(
        @{ Sheet = "zh-cn"; Xlf = $zhXlf; HandbackDate = $zhHandbackDate },
        @{ Sheet = "de-de"; Xlf = $deXlf; HandbackDate = $deHandbackDate }
    )) {

    $ws = $wb.Worksheets.Item($langInfo.Sheet)

    # Widen the Status column (C) and the target/handback file columns (I, J)
    # so the longer generated values are readable.
    $ws.Columns.Item(3).ColumnWidth  = 29.15
    $ws.Columns.Item(9).ColumnWidth  = 39.15
    $ws.Columns.Item(10).ColumnWidth = 39.15

    foreach ($row in 2..3) {
        # Status column
        $statusCell = $ws.Cells.Item($row, 3)
        if ($statusCell.Value() -eq $oldStatus) {
            $statusCell.Value = $newStatus
        }

        # Latest Target File (I) -> hyperlink back to the source file, same
        # as column A's link.
        $targetCell = $ws.Cells.Item($row, 9)
        $ws.Hyperlinks.Add($targetCell, $mdUrl, "", "", $mdFile)

        # Latest Handback File (J) -> the generated xlf for this language.
        $ws.Cells.Item($row, 10).Value = $langInfo.Xlf

        # Latest Handback DateTime (K) -> when the handback report was generated.
        $ws.Cells.Item($row, 11).Value = $langInfo.HandbackDate
    }
}
